# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates DAMSLTag (column I) and DialogAct (column J)
# values for the rows whose annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 5;  DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 12; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 14; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 21; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 23; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 27; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 28; DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 32; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 33; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 34; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 43; DAMSLTag = "qy"; DialogAct = "Yes-No-Question" },
    @{ Row = 49; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 55; DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 66; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 68; DAMSLTag = "sv"; DialogAct = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
